# Apply cryptos list update (Thu Dec  7 05:55:55 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.985.25'
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").Value = '2.265.41'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.658'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +5.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '233.33'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +1.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.67'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +0.29%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.453'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +7.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0983'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -0.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.95'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +1.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '27.01'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +4.72%  '

$ws.Range("E13").Value = '  +2.51%  '

$ws.Range("D14").Value = '2.602.17'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.70'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.50%  '

$ws.Range("E16").Value = '  +5.40%  '

$ws.Range("E17").Value = '  +3.38%  '

$ws.Range("D18").Value = '2.264.85'
$ws.Range("E18").Value = '  -0.19%  '

$ws.Range("D19").Value = '43.918.60'
$ws.Range("E19").Value = '  +0.75%  '

$ws.Range("D20").Value = '0.0₃0988'
$ws.Range("E20").Value = '  +3.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.33'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +1.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.16'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +1.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.24'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +1.09%  '

$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("E25").Value = '  -1.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +4.08%  '

$ws.Range("E27").Value = '  +18.66%  '

$ws.Range("E28").Value = '  +0.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.43'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +9.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.11'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +1.72%  '

$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("E32").Value = '  +0.32%  '

$ws.Range("E33").Value = '  +4.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.05'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +7.92%  '

$ws.Range("E35").Value = '  -0.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.01'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -2.13%  '

$ws.Range("E37").Value = '  -2.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.47'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -3.97%  '

$ws.Range("E39").Value = '  -0.96%  '

$ws.Range("E40").Value = '  +4.12%  '

$ws.Range("E41").Value = '  +0.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.83'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +5.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000225'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +1.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.31'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +3.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '98.66'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +1.69%  '

$ws.Range("E46").Value = '  -0.76%  '

$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.39'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -1.76%  '

$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.19'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -1.05%  '

$ws.Range("E49").Value = '  +1.55%  '

$ws.Range("D50").Value = '1.453.94'
$ws.Range("E50").Value = '  -1.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.05'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -5.76%  '
